$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (sheet name in workbook.xml)
$ws.Name = "estadisticas - 2025-05-04T11095"

# Column D (prv) gets a new value of 35 for all data rows (2-14)
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 4).Value = 35
}

# Row 2: PARTIDO POPULAR
$ws.Cells.Item(2, 8).Value = "PARTIDO POPULAR"
$ws.Cells.Item(2, 9).Value = "210.774"
$ws.Cells.Item(2, 10).Value = 48.05
$ws.Cells.Item(2, 11).Value = 29.36
$ws.Cells.Item(2, 12).Value = 48.49
$ws.Cells.Item(2, 14).Value = 4

# Row 3: COALICION CANARIA
$ws.Cells.Item(3, 8).Value = "COALICION CANARIA"
$ws.Cells.Item(3, 9).Value = "113.075"
$ws.Cells.Item(3, 10).Value = 25.78
$ws.Cells.Item(3, 11).Value = 15.75
$ws.Cells.Item(3, 12).Value = 26.01
$ws.Cells.Item(3, 14).Value = 2

# Row 4: PARTIDO SOCIALISTA OBRERO ESPAÑOL (dropped " - PROGRESISTAS")
$ws.Cells.Item(4, 8).Value = "PARTIDO SOCIALISTA OBRERO ESPAÃ‘OL"
$ws.Cells.Item(4, 9).Value = "80.695"
$ws.Cells.Item(4, 10).Value = 18.399999999999999
$ws.Cells.Item(4, 11).Value = 11.24
$ws.Cells.Item(4, 12).Value = 18.559999999999999
$ws.Cells.Item(4, 14).Value = 1

# Row 5: IZQUIERDA UNIDA CANARIA
$ws.Cells.Item(5, 8).Value = "IZQUIERDA UNIDA CANARIA"
$ws.Cells.Item(5, 9).Value = "10.941"
$ws.Cells.Item(5, 10).Value = 2.4900000000000002
$ws.Cells.Item(5, 11).Value = 1.52
$ws.Cells.Item(5, 12).Value = 2.52

# Row 6: PARTIDO DE INDEPENDIENTES DE LANZAROTE (votos unchanged at 10.323)
$ws.Cells.Item(6, 8).Value = "PARTIDO DE INDEPENDIENTES DE LANZAROTE"
$ws.Cells.Item(6, 9).Value = "10.323"
$ws.Cells.Item(6, 10).Value = 2.35
$ws.Cells.Item(6, 11).Value = 1.44
$ws.Cells.Item(6, 12).Value = 2.37

# Row 7: LOS VERDES DE CANARIAS
$ws.Cells.Item(7, 8).Value = "LOS VERDES DE CANARIAS"
$ws.Cells.Item(7, 9).Value = "4.982"
$ws.Cells.Item(7, 10).Value = 1.1399999999999999
$ws.Cells.Item(7, 11).Value = 0.69
$ws.Cells.Item(7, 12).Value = 1.1499999999999999

# Row 8: UNIDAD PROGRESISTAS DE CANARIAS (new party inserted here, votos now a number)
$ws.Cells.Item(8, 8).Value = "UNIDAD PROGRESISTAS DE CANARIAS"
$ws.Cells.Item(8, 9).Value = 980
$ws.Cells.Item(8, 10).Value = 0.22
$ws.Cells.Item(8, 11).Value = 0.14000000000000001
$ws.Cells.Item(8, 12).Value = 0.23

# Row 9: PARTIDO COMUNISTA DEL PUEBLO CANARIO
$ws.Cells.Item(9, 8).Value = "PARTIDO COMUNISTA DEL PUEBLO CANARIO"
$ws.Cells.Item(9, 9).Value = 752
$ws.Cells.Item(9, 10).Value = 0.17
$ws.Cells.Item(9, 11).Value = 0.1
$ws.Cells.Item(9, 12).Value = 0.17

# Row 10: PARTIDO HUMANISTA
$ws.Cells.Item(10, 8).Value = "PARTIDO HUMANISTA"
$ws.Cells.Item(10, 9).Value = 630
$ws.Cells.Item(10, 10).Value = 0.14000000000000001
$ws.Cells.Item(10, 11).Value = 0.09
$ws.Cells.Item(10, 12).Value = 0.14000000000000001

# Row 11: LA FALANGE
$ws.Cells.Item(11, 8).Value = "LA FALANGE"
$ws.Cells.Item(11, 9).Value = 478
$ws.Cells.Item(11, 10).Value = 0.11
$ws.Cells.Item(11, 11).Value = 0.07
$ws.Cells.Item(11, 12).Value = 0.11

# Row 12: UNION CENTRISTA-CENTRO DEMOCRATICO Y SOCIAL
$ws.Cells.Item(12, 8).Value = "UNION CENTRISTA-CENTRO DEMOCRATICO Y SOCIAL"
$ws.Cells.Item(12, 9).Value = 437
$ws.Cells.Item(12, 10).Value = 0.1
$ws.Cells.Item(12, 11).Value = 0.06
$ws.Cells.Item(12, 12).Value = 0.1

# Row 13: TAGOROR PENSIONISTA DE CANARIAS (votos unchanged at 319)
$ws.Cells.Item(13, 8).Value = "TAGOROR PENSIONISTA DE CANARIAS"
$ws.Cells.Item(13, 9).Value = 319
$ws.Cells.Item(13, 10).Value = 0.07
$ws.Cells.Item(13, 11).Value = 0.04
$ws.Cells.Item(13, 12).Value = 0.07

# Row 14: PARTIDO DE LA LEY NATURAL (votos unchanged at 301)
$ws.Cells.Item(14, 8).Value = "PARTIDO DE LA LEY NATURAL"
$ws.Cells.Item(14, 9).Value = 301
$ws.Cells.Item(14, 10).Value = 0.07
$ws.Cells.Item(14, 11).Value = 0.04
$ws.Cells.Item(14, 12).Value = 0.07
